$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like a plain number (e.g. "1.011").
# Excel auto-converts such text to a numeric Value on assignment, so the whole
# price column is temporarily switched to Text format, values are written, and
# then the format is restored to Normal/General so the cells keep their original
# (unstyled) look -- only the literal text content changes, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.969.00"
$ws.Range("D3").Value = "1.847.02"
$ws.Range("D5").Value = "1.011"
$ws.Range("D6").Value = "309.81"
$ws.Range("D7").Value = "0.4780"
$ws.Range("D9").Value = "0.07229"
$ws.Range("D10").Value = "0.9287"
$ws.Range("D11").Value = "19.70"
$ws.Range("D12").Value = "0.07721"
$ws.Range("D13").Value = "1.875.09"
$ws.Range("D15").Value = "6.430"
$ws.Range("D16").Value = "88.65"
$ws.Range("D17").Value = "1.013"
$ws.Range("D20").Value = "26.998.58"
$ws.Range("D21").Value = "14.47"
$ws.Range("D22").Value = "5.061"
$ws.Range("D24").Value = "1.931"
$ws.Range("D25").Value = "152.68"
$ws.Range("D26").Value = "18.21"
$ws.Range("D27").Value = "2.001"
$ws.Range("D29").Value = "4.974"
$ws.Range("D30").Value = "0.08896"
$ws.Range("D31").Value = "3.310"
$ws.Range("D32").Value = "1.177"
$ws.Range("D33").Value = "0.7408"
$ws.Range("D34").Value = "4.501"
$ws.Range("D35").Value = "2.740"
$ws.Range("D36").Value = "1.116"
$ws.Range("D37").Value = "0.01956"
$ws.Range("D38").Value = "0.05263"
$ws.Range("D40").Value = "0.5219"
$ws.Range("D41").Value = "6.986"
$ws.Range("D42").Value = "0.1513"
$ws.Range("D43").Value = "8.204"
$ws.Range("D44").Value = "10.56"
$ws.Range("D45").Value = "0.4751"
$ws.Range("D46").Value = "1.012"
$ws.Range("D47").Value = "101.65"
$ws.Range("D48").Value = "1.607"
$ws.Range("D49").Value = "65.60"
$ws.Range("D50").Value = "0.06058"
$ws.Range("D51").Value = "0.8877"

$ws.Range("D2:D51").Style = "Normal"

# Column E (percentage change) already carries padding spaces and a "%" sign,
# so Excel keeps it as text automatically -- no format juggling required.
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +2.53%  "
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  +4.00%  "

Write-Output "cryptos list updated"
